{"js": "// Iteration 5 -> Iteration 6: update the heading number and let the\n// \"_GoBack\" last-edit bookmark follow the edit (Word moves it to the\n// most-recently-edited location), removing it from its old spot.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document title heading (\"Iteration 5: Wrapping up\") is the first\n// paragraph (Heading 1 style). Find the standalone \"5\" run within it.\nconst headingPara = paragraphs.items[0];\nconst hits = headingPara.search(\"5\", { matchCase: true, matchWholeWord: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"5\" in the heading paragraph.');\n}\n\nconst target = hits.items[0];\ntarget.insertText(\"6\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Word keeps a single \"_GoBack\" bookmark marking the last edit; remove any\n// existing one before re-adding it at the new edit location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst endOfEdit = target.getRange(Word.RangeLocation.end);\nendOfEdit.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Iteration 5 -> Iteration 6: bump the heading number and let the\n# \"_GoBack\" last-edit bookmark follow the edit (Word keeps only one such\n# bookmark, relocating it to wherever was most recently edited).\n\n$d = $word.ActiveDocument\n\n# The document title (\"Iteration 5: Wrapping up\") is the first paragraph.\n$headingPara = $d.Paragraphs(1)\n$rng = $headingPara.Range\n$find = $rng.Find\n$find.Text = \"5\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    $editStart = $rng.Start\n    $editEnd = $rng.End\n\n    # Drop a transient bookmark right before the \"5\" so Word is forced to\n    # keep \"Iteration \" as its own run once we rewrite the text that\n    # follows it.\n    $preSplit = $d.Range($editStart, $editStart)\n    $d.Bookmarks.Add(\"ZZTMPSPLIT\", $preSplit)\n\n    # \"5\" -> \"6\"\n    $editRange = $d.Range($editStart, $editEnd)\n    $editRange.Text = \"6\"\n\n    # Word only ever keeps a single \"_GoBack\" bookmark (the most recent\n    # edit location); remove any existing one, then re-add it right after\n    # the text we just changed.\n    try {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    } catch {\n    }\n    $goBackRange = $d.Range($editStart + 1, $editStart + 1)\n    $d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n\n    # Clean up the transient helper bookmark.\n    $d.Bookmarks(\"ZZTMPSPLIT\").Delete()\n}\n"}
